$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two rows below are section-header rows that only hold a label in
# column A (no data in B/C/D):
#   row 5 -> "situação do domicílio"
#   row 8 -> "grandes regiões e unidades da federação"
# Both are being dropped; deleting them shifts every subsequent row up,
# which re-aligns the numeric data with the correct region/state labels
# and compacts the now-unused shared strings.
# Delete the higher-numbered row first so row 5's index isn't affected.
$ws.Rows(8).Delete()
$ws.Rows(5).Delete()
